$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (copy format from H1, which carries the bold/border
# header style).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Per-row data for the new I0 / IF columns.
$data = @(
    @(2, 7, 7),
    @(3, 7, 7),
    @(4, 7, 7),
    @(5, 6, 6),
    @(6, 7, 7),
    @(7, 8, 8),
    @(8, 9, 9),
    @(9, 8, 8),
    @(10, 7, 8),
    @(11, 7, 7),
    @(12, 7, 7),
    @(13, 6, 7),
    @(14, 7, 8),
    @(15, 6, 6),
    @(16, 8, 8),
    @(17, 5, 6),
    @(18, 6, 6),
    @(19, 7, 7),
    @(20, 8, 8),
    @(21, 5, 5),
    @(22, 8, 8),
    @(23, 6, 7),
    @(24, 7, 7),
    @(25, 7, 7),
    @(26, 7, 7),
    @(27, 7, 7),
    @(28, 6, 7),
    @(29, 6, 7),
    @(30, 6, 6),
    @(31, 7, 7),
    @(32, 6, 6),
    @(33, 6, 6),
    @(34, 4, 5),
    @(35, 6, 6),
    @(36, 9, 9),
    @(37, 7, 7),
    @(38, 8, 8),
    @(39, 7, 8),
    @(40, 7, 7),
    @(41, 5, 5),
    @(42, 7, 7),
    @(43, 5, 5),
    @(44, 5, 5),
    @(45, 6, 7),
    @(46, 4, 5),
    @(47, 5, 5),
    @(48, 5, 6),
    @(49, 7, 7),
    @(50, 7, 7),
    @(51, 7, 7),
    @(52, 8, 8),
    @(53, 8, 8),
    @(54, 7, 7),
    @(55, 6, 6),
    @(56, 7, 7),
    @(57, 6, 6),
    @(58, 7, 8),
    @(59, 9, 9),
    @(60, 6, 6),
    @(61, 9, 9),
    @(62, 7, 7),
    @(63, 1, 1),
    @(64, 6, 6),
    @(65, 3, 3),
    @(66, 7, 7),
    @(67, 4, 4),
    @(68, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
